# Insert a new weekly price record as row 83 on the single data sheet.
# This shifts the existing rows 83..216 down to 84..217 (dimension grows
# from A1:R216 to A1:R217), matching the authors' "Fruta / hortaliza,
# semanal" commit that adds one more weekly observation to the series.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push row 83 (and everything below it) down by one row.
$ws.Rows.Item(83).Insert()

# Populate the freshly inserted row 83 with the new observation.
$ws.Cells.Item(83, 1).Value  = 10
$ws.Cells.Item(83, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(83, 3).Value  = "La Araucanía"
$ws.Cells.Item(83, 4).Value  = 44894
$ws.Cells.Item(83, 5).Value  = 9
$ws.Cells.Item(83, 6).Value  = 100114007
$ws.Cells.Item(83, 7).Value  = "Jengibre"
$ws.Cells.Item(83, 8).Value  = "Sin especificar"
$ws.Cells.Item(83, 9).Value  = "Primera"
$ws.Cells.Item(83, 10).Value = 20
$ws.Cells.Item(83, 11).Value = 20000
$ws.Cells.Item(83, 12).Value = 20000
$ws.Cells.Item(83, 13).Value = 20000
$ws.Cells.Item(83, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item(83, 15).Value = "Perú"
$ws.Cells.Item(83, 16).Value = 1538
$ws.Cells.Item(83, 17).Value = 13
$ws.Cells.Item(83, 18).Value = "Hortaliza"
